# Update the "修改时间" (modified time) timestamps across the three
# portfolio sheets from 202509211559 to 202509211606.
#
# The timestamps are stored as text (not numbers), so the all-digit
# string is written with a leading apostrophe - exactly as a user would
# type it into Excel - to force text interpretation instead of letting
# the value be auto-coerced into a number.

$wb = $excel.ActiveWorkbook

$newTimestamp = "'202509211606"

# Sheet 1: 大智投资组合 - timestamps in column E, rows 2-9
$ws1 = $wb.Worksheets.Item("大智投资组合")
for ($r = 2; $r -le 9; $r++) {
    $ws1.Cells.Item($r, 5).Value = $newTimestamp
}

# Sheet 2: 大成投资组合 - timestamps in column E, rows 2-11
$ws2 = $wb.Worksheets.Item("大成投资组合")
for ($r = 2; $r -le 11; $r++) {
    $ws2.Cells.Item($r, 5).Value = $newTimestamp
}

# Sheet 3: 我的投资组合 - timestamps in column G, rows 2-13
$ws3 = $wb.Worksheets.Item("我的投资组合")
for ($r = 2; $r -le 13; $r++) {
    $ws3.Cells.Item($r, 7).Value = $newTimestamp
}
